# Insert a new row at the top of the data (row 2), pushing all existing
# price rows down by one, and fill it in with today's date + the latest
# price figures (carried forward, same as the rest of the series).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Insert()
$ws.Range("A2:D2").ClearFormats()

$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610

# Writing the date string straight into .Value lets Excel's text parser
# recognize it as a real date and auto-convert the cell, which is not what
# the source data uses (dates are plain text here). Route the literal text
# through a TRIM() formula + paste-values round trip so the cell keeps a
# plain-text type with no special number formatting, matching the rest of
# the column.
$ws.Range("A2").Value = "2025-12-22 "
$ws.Range("F2").Formula = "=TRIM(A2)"
$ws.Range("F2").Copy()
$ws.Range("A2").PasteSpecial(-4163)
$ws.Range("F2").Clear()
